# The deck ships with two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme"  (targeted by the Notes Master)
#   ppt/theme/theme2.xml -> "Integral"      (targeted by the Presentation / Slide Master)
#
# The authored change swaps the two themes' content: the theme that backs
# the Slide Master becomes the plain "Office Theme" palette, while the
# (COM-unreachable) Notes Master theme part becomes "Integral". The only
# part of that swap that is reachable from the PowerPoint object model is
# the live/active theme color scheme hanging off the Slide Master (it is
# the single ThemeColorScheme the host exposes), so recolor it here to the
# stock Office Theme palette.

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$scheme = $design.SlideMaster.Theme.ThemeColorScheme

$scheme.Item(1).RGB  = 0          # dk1      000000
$scheme.Item(2).RGB  = 16777215   # lt1      FFFFFF
$scheme.Item(3).RGB  = 6968388    # dk2      44546A
$scheme.Item(4).RGB  = 15132391   # lt2      E7E6E6
$scheme.Item(5).RGB  = 13998939   # accent1  5B9BD5
$scheme.Item(6).RGB  = 3243501    # accent2  ED7D31
$scheme.Item(7).RGB  = 10855845   # accent3  A5A5A5
$scheme.Item(8).RGB  = 49407      # accent4  FFC000
$scheme.Item(9).RGB  = 12874308   # accent5  4472C4
$scheme.Item(10).RGB = 4697456    # accent6  70AD47
$scheme.Item(11).RGB = 12673797   # hlink    0563C1
$scheme.Item(12).RGB = 7491477    # folHlink 954F72
